$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44294
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 12000
$ws.Cells.Item(2, 15).Value = 12000
$ws.Cells.Item(2, 16).Value = 12000
$ws.Cells.Item(2, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(2, 19).Value = 800

# Row 3
$ws.Cells.Item(3, 4).Value = 44348
$ws.Cells.Item(3, 13).Value = 200

# Row 4
$ws.Cells.Item(4, 4).Value = 44340
$ws.Cells.Item(4, 13).Value = 230
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 15).Value = 20000
$ws.Cells.Item(4, 16).Value = 20000
$ws.Cells.Item(4, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(4, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(4, 19).Value = 1111
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = 44354
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 100

# Row 6
$ws.Cells.Item(6, 4).Value = 44691
$ws.Cells.Item(6, 14).Value = 17000
$ws.Cells.Item(6, 15).Value = 17000
$ws.Cells.Item(6, 16).Value = 17000
$ws.Cells.Item(6, 19).Value = 944

# Row 7
$ws.Cells.Item(7, 4).Value = 44358
$ws.Cells.Item(7, 14).Value = 18000
$ws.Cells.Item(7, 15).Value = 18000
$ws.Cells.Item(7, 16).Value = 18000
$ws.Cells.Item(7, 19).Value = 1000

# Row 8
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 17000
$ws.Cells.Item(8, 15).Value = 17000
$ws.Cells.Item(8, 16).Value = 17000
$ws.Cells.Item(8, 19).Value = 944

# Row 9
$ws.Cells.Item(9, 4).Value = 44326
$ws.Cells.Item(9, 12).Value = 'Especial'
$ws.Cells.Item(9, 13).Value = 300
$ws.Cells.Item(9, 14).Value = 20000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 20000
$ws.Cells.Item(9, 19).Value = 1111

# Row 10
$ws.Cells.Item(10, 4).Value = 44316
$ws.Cells.Item(10, 12).Value = 'Especial'
$ws.Cells.Item(10, 13).Value = 300

# Row 11
$ws.Cells.Item(11, 4).Value = 44692
$ws.Cells.Item(11, 12).Value = 'Especial'
$ws.Cells.Item(11, 13).Value = 150
$ws.Cells.Item(11, 14).Value = 17000
$ws.Cells.Item(11, 15).Value = 17000
$ws.Cells.Item(11, 16).Value = 17000
$ws.Cells.Item(11, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(11, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 19).Value = 944
$ws.Cells.Item(11, 20).Value = 18

# Row 12
$ws.Cells.Item(12, 4).Value = 44355
$ws.Cells.Item(12, 13).Value = 50
$ws.Cells.Item(12, 14).Value = 18000
$ws.Cells.Item(12, 15).Value = 18000
$ws.Cells.Item(12, 16).Value = 18000
$ws.Cells.Item(12, 19).Value = 1000

# Row 13
$ws.Cells.Item(13, 4).Value = 44680
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 14).Value = 15000
$ws.Cells.Item(13, 15).Value = 15000
$ws.Cells.Item(13, 16).Value = 15000
$ws.Cells.Item(13, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(13, 19).Value = 1000
$ws.Cells.Item(13, 20).Value = 15

# Row 14
$ws.Cells.Item(14, 4).Value = 44328
$ws.Cells.Item(14, 12).Value = 'Especial'
$ws.Cells.Item(14, 13).Value = 250

# Row 15
$ws.Cells.Item(15, 4).Value = 44319
$ws.Cells.Item(15, 13).Value = 120

# Row 16
$ws.Cells.Item(16, 4).Value = 44342

# Row 17
$ws.Cells.Item(17, 4).Value = 44714
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 20000
$ws.Cells.Item(17, 15).Value = 20000
$ws.Cells.Item(17, 16).Value = 20000
$ws.Cells.Item(17, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(17, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(17, 19).Value = 1111
$ws.Cells.Item(17, 20).Value = 18

# Row 18
$ws.Cells.Item(18, 4).Value = 44299
$ws.Cells.Item(18, 14).Value = 15000
$ws.Cells.Item(18, 15).Value = 15000
$ws.Cells.Item(18, 16).Value = 15000
$ws.Cells.Item(18, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(18, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(18, 19).Value = 1000
$ws.Cells.Item(18, 20).Value = 15

# Row 19
$ws.Cells.Item(19, 4).Value = 44291
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 150
$ws.Cells.Item(19, 14).Value = 12000
$ws.Cells.Item(19, 15).Value = 12000
$ws.Cells.Item(19, 16).Value = 12000
$ws.Cells.Item(19, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(19, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 19).Value = 800
$ws.Cells.Item(19, 20).Value = 15
